# Actualización al 8 de julio de 2023
$wb = $excel.ActiveWorkbook

$wsIngreso = $wb.Worksheets.Item("Ingreso")
$wsGastos  = $wb.Worksheets.Item("Gastos")

# --- Sheet "Ingreso": append new contribution rows for 2023-07-08 (serial 45115) ---
$fecha = 45115

$ingresoRows = @(
    @("Michy",      200),
    @("Anuel",      100),
    @("Mac Daniel",  80),
    @("Carlos",     100),
    @("Invitados",  100),
    @("Chamo",      200),
    @("Randy",      100),
    @("Jeremy",     100),
    @("Johan",      300)
)

$startRow = 474
for ($i = 0; $i -lt $ingresoRows.Count; $i++) {
    $row = $startRow + $i
    $nombre = $ingresoRows[$i][0]
    $monto  = $ingresoRows[$i][1]

    $wsIngreso.Cells.Item($row, 1).Value = $fecha
    $wsIngreso.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $wsIngreso.Cells.Item($row, 2).Value = $nombre
    $wsIngreso.Cells.Item($row, 3).Value = $monto
    $wsIngreso.Cells.Item($row, 4).Value = "Aporte"
}

$lastRow = $startRow + $ingresoRows.Count - 1
$wsIngreso.Activate()
$wsIngreso.Range("A$lastRow").Select()

# --- Sheet "Gastos": fix date on row 49 and append new referee/water/ice expense ---
$wsGastos.Cells.Item(49, 1).Value = 45109

$wsGastos.Cells.Item(50, 1).Value = 45115
$wsGastos.Cells.Item(50, 1).NumberFormat = "yyyy\-mm\-dd;@"
$wsGastos.Cells.Item(50, 2).Value = "Arbitro, agua y hielo"
$wsGastos.Cells.Item(50, 3).Value = 941
